$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main_Loop_Parameters")

# --- Simulation size updates ---
$ws.Range("B4").Value = 10000    # n_households: 100000 -> 10000
$ws.Range("B5").Value = 450      # n_consumer_firms: 4500 -> 450
$ws.Range("B6").Value = 50       # n_capital_firms: 500 -> 50

# --- Insert two new Dynamic Bank Parameters rows after
#     "bank_leverage_ratio_upper_threshold" (row 44) ---
$ws.Range("A45:A46").EntireRow.Insert()

$ws.Range("A45").Value = "bank_max_interest_rate"
$ws.Range("B45").Value = 1

$ws.Range("A46").Value = "bank_max_interest_rate_change"
$ws.Range("B46").Value = 0.05

# Match row height/formatting of the surrounding data rows
$ws.Range("A45:B46").RowHeight = 19.5

Write-Output "edit complete"
